$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 81.295946
$ws.Range("H2").Value = 162.591892
$ws.Range("I2").Value = 0.2087576831039754
$ws.Range("J2").Value = 0.1513633345147276
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.032708
$ws.Range("N2").Value = 0.065416
$ws.Range("Q2").Value = 2.659027801768
$ws.Range("R2").Value = 10.636111207072
$ws.Range("S2").Value = 0.2087576831039754
$ws.Range("T2").Value = 0.1513633345147276

# Row 3
$ws.Range("I3").Value = 0.08203258389144356
$ws.Range("J3").Value = 0.08921869546581183
$ws.Range("K3").Value = 2
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 0.032708
$ws.Range("N3").Value = 0.065416
$ws.Range("Q3").Value = 1.044880925937333
$ws.Range("R3").Value = 6.269285555624
$ws.Range("S3").Value = 0.08203258389144356
$ws.Range("T3").Value = 0.08921869546581183

# Row 4
$ws.Range("G4").Value = 52.98516133333333
$ws.Range("H4").Value = 158.955484
$ws.Range("I4").Value = 0.1360591771554883
$ws.Range("J4").Value = 0.1479780559884402
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 0.032708
$ws.Range("N4").Value = 0.065416
$ws.Range("Q4").Value = 1.733038656890667
$ws.Range("R4").Value = 10.398231941344
$ws.Range("S4").Value = 0.1360591771554883
$ws.Range("T4").Value = 0.1479780559884402

# Row 5
$ws.Range("G5").Value = 12.8032175
$ws.Range("H5").Value = 25.606435
$ws.Range("I5").Value = 0.03287703942305158
$ws.Range("J5").Value = 0.02383806067423477
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.032708
$ws.Range("N5").Value = 0.065416
$ws.Range("Q5").Value = 0.41876763799
$ws.Range("R5").Value = 1.67507055196
$ws.Range("S5").Value = 0.03287703942305158
$ws.Range("T5").Value = 0.02383806067423477

# Row 6
$ws.Range("G6").Value = 30.251696
$ws.Range("H6").Value = 90.755088
$ws.Range("I6").Value = 0.07768252019511281
$ws.Range("J6").Value = 0.08448756315510209
$ws.Range("K6").Value = 2
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 0.032708
$ws.Range("N6").Value = 0.065416
$ws.Range("Q6").Value = 0.989472472768
$ws.Range("R6").Value = 5.936834836608
$ws.Range("S6").Value = 0.07768252019511281
$ws.Range("T6").Value = 0.08448756315510209

# Row 7
$ws.Range("G7").Value = 180.145574
$ws.Range("H7").Value = 540.436722
$ws.Range("I7").Value = 0.4625909962309284
$ws.Range("J7").Value = 0.5031142902016837
$ws.Range("K7").Value = 2
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 0.032708
$ws.Range("N7").Value = 0.065416
$ws.Range("Q7").Value = 5.892201434392001
$ws.Range("R7").Value = 35.35320860635201
$ws.Range("S7").Value = 0.4625909962309284
$ws.Range("T7").Value = 0.5031142902016837
